$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 with new values (per Dr Hou's revised analysis) ---

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.666083666666667
$ws.Range("H2").Value = 4.998251
$ws.Range("I2").Value = 0.6125276070882968
$ws.Range("J2").Value = 0.6125276070882968
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.55727433333333
$ws.Range("N2").Value = 34.671823
$ws.Range("O2").Value = 0.5239815261112396
$ws.Range("P2").Value = 0.5239815261112395
$ws.Range("Q2").Value = 19.25538599795255
$ws.Range("R2").Value = 173.298473981573
$ws.Range("S2").Value = 0.3209531503473915
$ws.Range("T2").Value = 0.3209531503473914

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.666083666666667
$ws.Range("H3").Value = 4.998251
$ws.Range("I3").Value = 0.6125276070882968
$ws.Range("J3").Value = 0.6125276070882968
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.09477133333333332
$ws.Range("N3").Value = 0.284314
$ws.Range("O3").Value = 0.004296724853919303
$ws.Range("P3").Value = 0.004296724853919302
$ws.Range("Q3").Value = 0.1578969705348889
$ws.Range("R3").Value = 1.421072734814
$ws.Range("S3").Value = 0.002631862593088002
$ws.Range("T3").Value = 0.002631862593088002

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.666083666666667
$ws.Range("H4").Value = 4.998251
$ws.Range("I4").Value = 0.6125276070882968
$ws.Range("J4").Value = 0.6125276070882968
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.404599
$ws.Range("N4").Value = 31.213797
$ws.Range("O4").Value = 0.471721749034841
$ws.Range("P4").Value = 0.471721749034841
$ws.Range("Q4").Value = 17.33493245211633
$ws.Range("R4").Value = 156.014392069047
$ws.Range("S4").Value = 0.2889425941478173
$ws.Range("T4").Value = 0.2889425941478173

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.053930333333333
$ws.Range("H5").Value = 3.161791
$ws.Range("I5").Value = 0.3874723929117032
$ws.Range("J5").Value = 0.3874723929117031
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.55727433333333
$ws.Range("N5").Value = 34.671823
$ws.Range("O5").Value = 0.5239815261112396
$ws.Range("P5").Value = 0.5239815261112395
$ws.Range("Q5").Value = 12.18056199055478
$ws.Range("R5").Value = 109.625057914993
$ws.Range("S5").Value = 0.2030283757638481
$ws.Range("T5").Value = 0.203028375763848

# --- Add new rows 6-7 (additional target-cluster combinations) ---

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.053930333333333
$ws.Range("H6").Value = 3.161791
$ws.Range("I6").Value = 0.3874723929117032
$ws.Range("J6").Value = 0.3874723929117031
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09477133333333332
$ws.Range("N6").Value = 0.284314
$ws.Range("O6").Value = 0.004296724853919303
$ws.Range("P6").Value = 0.004296724853919302
$ws.Range("Q6").Value = 0.09988238293044444
$ws.Range("R6").Value = 0.8989414463739999
$ws.Range("S6").Value = 0.001664862260831301
$ws.Range("T6").Value = 0.0016648622608313

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.053930333333333
$ws.Range("H7").Value = 3.161791
$ws.Range("I7").Value = 0.3874723929117032
$ws.Range("J7").Value = 0.3874723929117031
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.404599
$ws.Range("N7").Value = 31.213797
$ws.Range("O7").Value = 0.471721749034841
$ws.Range("P7").Value = 0.471721749034841
$ws.Range("Q7").Value = 10.96572249226967
$ws.Range("R7").Value = 98.69150243042699
$ws.Range("S7").Value = 0.1827791548870238
$ws.Range("T7").Value = 0.1827791548870237
